# Regenerate merged AHB files
#
# 1. Rename the "_old"/"_new" header-suffix columns to "_FV2410"/"_FV2504"
# 2. Freeze the header row (pane split after row 1)
# 3. Convert the data range A1:U71 into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels -------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value2
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2410")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2504")
    }
}

# --- 2. Freeze header row -----------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the range into a Table -------------------------------------------
$tableRange = $ws.Range("A1:U71")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"
